# Charlie_User_Stories.docx edit script
#
# Summary of the change (from the commit message / diff):
#  - Highlight (green) the "(5 points) ... tournament logistics ... download"
#    user-story paragraph, but split the trailing space off of the final run
#    so that single trailing space stays un-highlighted.
#  - Highlight (green) the "(10 points) ... location of a tournament on a
#    map ... Google Maps API" user-story paragraph.
#  - Move the "_GoBack" bookmark (an empty/zero-length bookmark) from its
#    old spot in the Tournament Director "...format, and |allow teams..."
#    paragraph to the very end of the newly-highlighted Google Maps
#    paragraph (right after "Google Maps API", before the paragraph mark).

$d = $word.ActiveDocument

function Find-ParagraphByText($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text -like "*$needle*") {
            return $p
        }
    }
    return $null
}

# wdColorIndex value that this host serialises as <w:highlight w:val="green"/>
$GREEN = 4
$NOHIGHLIGHT = 0

# ---------------------------------------------------------------------
# 1) "(5 points) ... tournament logistics ..." paragraph
# ---------------------------------------------------------------------
$pLogistics = Find-ParagraphByText $d "tournament logistics"

# Highlight the whole paragraph green first (covers all 4 existing runs).
$pLogistics.Range.HighlightColorIndex = $GREEN

# Now split off the trailing space of "... format that I can download "
# into its own, non-highlighted run/char, matching the target XML where
# the highlighted run ends with "...download" (no trailing space) and a
# plain trailing-space run follows it before the paragraph mark.
$pEnd = $pLogistics.Range.End
$trailingSpace = $d.Range($pEnd - 2, $pEnd - 1)
if ($trailingSpace.Text -eq " ") {
    $trailingSpace.HighlightColorIndex = $NOHIGHLIGHT
}

# ---------------------------------------------------------------------
# 2) "(10 points) ... location of a tournament on a map ... Google Maps
#    API" paragraph
# ---------------------------------------------------------------------
$pMap = Find-ParagraphByText $d "location of a tournament on a map"
$pMap.Range.HighlightColorIndex = $GREEN

# ---------------------------------------------------------------------
# 3) Move the "_GoBack" bookmark from the Tournament Director paragraph
#    to the end of the Google Maps paragraph (an empty bookmark).
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $oldBm = $d.Bookmarks.Item("_GoBack")
    $oldBm.Delete()
}

# Re-resolve the Google Maps paragraph range (unaffected in length by the
# highlight-only edits above) and insert a temporary one-character marker
# right before its paragraph mark so we can anchor a bookmark there, then
# delete the marker text -- this leaves a proper zero-length bookmark at
# that position (inserting a bookmark directly over a zero-length Range
# at a paragraph-mark boundary is not reliable on this host).
$pMap = Find-ParagraphByText $d "location of a tournament on a map"
$mapEnd = $pMap.Range.End
$insertPoint = $d.Range($mapEnd - 1, $mapEnd - 1)
$insertPoint.InsertAfter("X")

$marker = $d.Range($mapEnd - 1, $mapEnd)
$d.Bookmarks.Add("_GoBack", $marker)
$newBm = $d.Bookmarks.Item("_GoBack")
$newBm.Range.Text = ""
